$p = $ppt.ActivePresentation

function Set-ReadingText($shape, $newText) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $len = $full.Length
    $lead = 0
    while ($lead -lt $len -and $full[$lead] -eq "`r") { $lead++ }
    $rest = $len - $lead
    if ($rest -le 0) {
        $tr.Text = $newText
    } else {
        $sub = $tr.Characters($lead + 1, $rest)
        $sub.Text = $newText
    }
}

# Slide 1
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "殺す"
Set-ReadingText $s.Shapes.Item(2) "ころす"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "to kill, to slay, to murder, to slaughter | to suppress, to block, to hamper, to destroy (e.g. talent), to eliminate (e.g..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 2
$s = $p.Slides.Item(2)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "殺人"
Set-ReadingText $s.Shapes.Item(2) "さつじん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "murder, homicide, manslaughter..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 3
$s = $p.Slides.Item(3)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "農薬"
Set-ReadingText $s.Shapes.Item(2) "のうやく"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "agricultural chemical (i.e. pesticide, herbicide, fungicide, etc.), agrochemical, agrichemical..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 4
$s = $p.Slides.Item(4)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "収入印紙"
Set-ReadingText $s.Shapes.Item(2) "しゅうにゅういんし"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "revenue stamp..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 5
$s = $p.Slides.Item(5)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "収穫"
Set-ReadingText $s.Shapes.Item(2) "しゅうかく"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "harvest, crop, ingathering | fruits (of one's labors), gain, result, returns..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 6
$s = $p.Slides.Item(6)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "少量"
Set-ReadingText $s.Shapes.Item(2) "しょうりょう"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "small quantity, small amount | narrowmindedness..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 7
$s = $p.Slides.Item(7)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "完全"
Set-ReadingText $s.Shapes.Item(2) "かんぜん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "perfect, complete..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 8
$s = $p.Slides.Item(8)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "原因"
Set-ReadingText $s.Shapes.Item(2) "げんいん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "cause, origin, source..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 9
$s = $p.Slides.Item(9)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "一環"
Set-ReadingText $s.Shapes.Item(2) "いっかん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "link (e.g. in a chain of events), part (of a plan, campaign, activities, etc.) | monocyclic..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "67-68"

# Slide 10
$s = $p.Slides.Item(10)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "幅広い"
Set-ReadingText $s.Shapes.Item(2) "はばひろい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "extensive, wide, broad..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 11
$s = $p.Slides.Item(11)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "増幅"
Set-ReadingText $s.Shapes.Item(2) "ぞうふく"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "amplification (elec.) | magnification, amplification, making larger..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 12
$s = $p.Slides.Item(12)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "支える"
Set-ReadingText $s.Shapes.Item(2) "ささえる"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "to support, to prop, to sustain, to underlay, to hold up, to defend | to hold at bay, to stem, to check..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 13
$s = $p.Slides.Item(13)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "支持"
Set-ReadingText $s.Shapes.Item(2) "しじ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "support, backing, endorsement, approval | propping up, holding up, support..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 14
$s = $p.Slides.Item(14)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "反対"
Set-ReadingText $s.Shapes.Item(2) "はんたい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "opposition, resistance, antagonism, hostility, objection, dissent | reverse, opposite, inverse, contrary..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 15
$s = $p.Slides.Item(15)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "論文"
Set-ReadingText $s.Shapes.Item(2) "ろんぶん"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "thesis, essay, treatise, paper, article..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 16
$s = $p.Slides.Item(16)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "否定"
Set-ReadingText $s.Shapes.Item(2) "ひてい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "denial, negation, repudiation, disavowal | negation | NOT operation..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 17
$s = $p.Slides.Item(17)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "野菜"
Set-ReadingText $s.Shapes.Item(2) "やさい"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "vegetable..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 18
$s = $p.Slides.Item(18)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "穀物"
Set-ReadingText $s.Shapes.Item(2) "こくもつ"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "grain, cereal, corn..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Slide 19
$s = $p.Slides.Item(19)
$s.Shapes.Item(1).TextFrame.TextRange.Text = "害虫"
Set-ReadingText $s.Shapes.Item(2) "がいちゅう"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "harmful insect, noxious insect, vermin, pest..."
$s.Shapes.Item(4).TextFrame.TextRange.Text = "65-66"

# Remove the last 4 slides (20-23), from highest index to lowest
$p.Slides.Item(23).Delete()
$p.Slides.Item(22).Delete()
$p.Slides.Item(21).Delete()
$p.Slides.Item(20).Delete()

Write-Output "Final slide count: $($p.Slides.Count)"